$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly records arrived; push the existing rows (old 381..392)
# down by two rows (to 383..394) and insert the two new records at the
# top of the range (new rows 381 and 382).
$ws.Rows(381).Insert()
$ws.Rows(381).Insert()

# New row 381
$ws.Cells.Item(381, 1).Value = 5
$ws.Cells.Item(381, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(381, 3).Value = "Maule"
$ws.Cells.Item(381, 4).Value = 44509
$ws.Cells.Item(381, 5).Value = 7
$ws.Cells.Item(381, 6).Value = "Fruta"
$ws.Cells.Item(381, 7).Value = 100104
$ws.Cells.Item(381, 8).Value = "Frutos de pepita"
$ws.Cells.Item(381, 9).Value = 100104005
$ws.Cells.Item(381, 10).Value = "Pera"
$ws.Cells.Item(381, 11).Value = "Packham's Triumph"
$ws.Cells.Item(381, 12).Value = "Especial"
$ws.Cells.Item(381, 13).Value = 300
$ws.Cells.Item(381, 14).Value = 12000
$ws.Cells.Item(381, 15).Value = 12000
$ws.Cells.Item(381, 16).Value = 12000
$ws.Cells.Item(381, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(381, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(381, 19).Value = 667
$ws.Cells.Item(381, 20).Value = 18

# New row 382
$ws.Cells.Item(382, 1).Value = 5
$ws.Cells.Item(382, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(382, 3).Value = "Maule"
$ws.Cells.Item(382, 4).Value = 44509
$ws.Cells.Item(382, 5).Value = 7
$ws.Cells.Item(382, 6).Value = "Fruta"
$ws.Cells.Item(382, 7).Value = 100104
$ws.Cells.Item(382, 8).Value = "Frutos de pepita"
$ws.Cells.Item(382, 9).Value = 100104005
$ws.Cells.Item(382, 10).Value = "Pera"
$ws.Cells.Item(382, 11).Value = "Packham's Triumph"
$ws.Cells.Item(382, 12).Value = "Especial"
$ws.Cells.Item(382, 13).Value = 200
$ws.Cells.Item(382, 14).Value = 12000
$ws.Cells.Item(382, 15).Value = 12000
$ws.Cells.Item(382, 16).Value = 12000
$ws.Cells.Item(382, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(382, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(382, 19).Value = 667
$ws.Cells.Item(382, 20).Value = 18
